$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database
# Applying updated values to rows 2-4 (Australia Semiconductor Equip companies)

$ws.Range("D2").Value = 0.147
$ws.Range("E2").Value = -0.0285
$ws.Range("G2").Value = -7.699779249448123
$ws.Range("H2").Value = -7.699779249448123
$ws.Range("I2").Value = -9.015452538631346
$ws.Range("J2").Value = -9.015452538631346
$ws.Range("K2").Value = -4.02
$ws.Range("L2").Value = -8.874172185430462
$ws.Range("U2").Value = 4.521
$ws.Range("V2").Value = 0.07721605465414175
$ws.Range("W2").Value = -0.2307757072294567
$ws.Range("X2").Value = 0.09169437868049965
$ws.Range("Y2").Value = -0.3224700859099563
$ws.Range("Z2").Value = 0.07389885807504078
$ws.Range("AA2").Value = -0.5335391369966385
$ws.Range("AB2").Value = 0.09086785821437925
$ws.Range("AC2").Value = -0.6244069952110177
$ws.Range("AD2").Value = 1.19
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1.19
$ws.Range("AG2").Value = -3.331
$ws.Range("AH2").Value = 0.01991965182457315
$ws.Range("AI2").Value = 0.08173076923076922
$ws.Range("AJ2").Value = -0.06032343939586012
$ws.Range("AK2").Value = -0.3318059567686024
$ws.Range("AL2").Value = 0.062
$ws.Range("AM2").Value = 0.037
$ws.Range("AN2").Value = -0.3411697247706422
$ws.Range("AO2").Value = -65.87096774193547
$ws.Range("AP2").Value = 0.9549885321100917
$ws.Range("AQ2").Value = -110.3783783783784
$ws.Range("E3").Value = -0.0285
$ws.Range("K3").Value = 0.12
$ws.Range("O3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 0.771
$ws.Range("V3").Value = 0.09233532934131737
$ws.Range("W3").Value = 0.04580152671755724
$ws.Range("X3").Value = 0.09101639306868477
$ws.Range("Y3").Value = -0.04521486635112752
$ws.Range("AA3").Value = -0.02831858407079646
$ws.Range("AB3").Value = 0.09101639306868477
$ws.Range("AC3").Value = -0.1193349771394812
$ws.Range("AG3").Value = -0.771
$ws.Range("AJ3").Value = -0.1017284602190263
$ws.Range("AK3").Value = -0.1908888338697698
$ws.Range("AM3").Value = -0.004
$ws.Range("AN3").Value = -0
$ws.Range("AP3").Value = 16.0625
$ws.Range("AQ3").Value = 16
$ws.Range("D4").Value = 0.147
$ws.Range("G4").Value = -7.593818984547461
$ws.Range("H4").Value = -7.593818984547461
$ws.Range("I4").Value = -8.874172185430462
$ws.Range("J4").Value = -8.874172185430462
$ws.Range("K4").Value = -4.14
$ws.Range("L4").Value = -9.139072847682119
$ws.Range("U4").Value = 3.75
$ws.Range("V4").Value = 0.0747011952191235
$ws.Range("W4").Value = -0.5073529411764706
$ws.Range("X4").Value = 0.09237236429231452
$ws.Range("Y4").Value = -0.5997253054687851
$ws.Range("Z4").Value = 0.1170542635658915
$ws.Range("AA4").Value = -1.03875968992248
$ws.Range("AB4").Value = 0.09071932336007373
$ws.Range("AC4").Value = -1.129479013282554
$ws.Range("AD4").Value = 1.19
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 1.19
$ws.Range("AG4").Value = -2.56
$ws.Range("AH4").Value = 0.0231562560809496
$ws.Range("AI4").Value = 0.122051282051282
$ws.Range("AJ4").Value = -0.05373635600335852
$ws.Range("AK4").Value = -0.4266666666666667
$ws.Range("AL4").Value = 0.062
$ws.Range("AM4").Value = 0.04099999999999999
$ws.Range("AN4").Value = -0.3459302325581395
$ws.Range("AO4").Value = -64.83870967741935
$ws.Range("AP4").Value = 0.7441860465116279
$ws.Range("AQ4").Value = -98.04878048780488
